# Applies the update: several rows in the existing table had their
# match-details columns (F:V) shuffled amongst sibling rows that share the
# same kickoff date/time, and three brand-new match rows (128-130) were
# appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($row1, $row2)
    $r1 = $ws.Range("F${row1}:V${row1}")
    $r2 = $ws.Range("F${row2}:V${row2}")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

function Rotate-Rows {
    # Each row in $rows (in order) receives the F:V content that currently
    # belongs to the next row in the list (wrapping around at the end).
    param($rows)
    $originals = @()
    foreach ($r in $rows) {
        $originals += , $ws.Range("F${r}:V${r}").Value2
    }
    $count = $rows.Length
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $rows[$i]
        $srcIndex = ($i + 1) % $count
        $ws.Range("F${destRow}:V${destRow}").Value2 = $originals[$srcIndex]
    }
}

# Simple pairwise swaps
Swap-Rows 9 10
Swap-Rows 61 62
Swap-Rows 81 82
Swap-Rows 83 84
Swap-Rows 85 86
Swap-Rows 92 93
Swap-Rows 95 96
Swap-Rows 98 99
Swap-Rows 119 120

# Three-way rotation: 124 <- 125 <- 126 <- 124
Rotate-Rows @(124, 125, 126)

# Append three new match rows at the bottom of the table (128, 129, 130),
# copying the formatting (styles/number formats) from the last existing
# data row (127) first.
$ws.Range("A127:V127").Copy()
$ws.Range("A128:V130").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @(127,"saudi-arabia","saudi-professional-league","2023-2024",45260.66666666666,"Al Riyadh",0,"Al Hazem",0,1.85,"28/11/2023 11:42",2.83,"30/11/2023 15:59",4.24,"28/11/2023 11:42",3.38,"30/11/2023 15:58",3.34,"28/11/2023 11:42",2.54,"30/11/2023 15:59","https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-riyadh-al-hazem-rass/YJiZWNQF/"),
    @(128,"saudi-arabia","saudi-professional-league","2023-2024",45260.66666666666,"Abha",0,"Al Ahli SC",6,4.13,"28/11/2023 11:42",5.16,"30/11/2023 15:59",4.2,"28/11/2023 11:42",4.65,"30/11/2023 15:59",1.67,"28/11/2023 11:42",1.56,"30/11/2023 15:59","https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/abha-al-ahli-sc/tQQIOaYe/"),
    @(129,"saudi-arabia","saudi-professional-league","2023-2024",45260.79166666666,"Al Ittihad",4,"Al Khaleej",2,1.23,"25/11/2023 19:13",1.29,"30/11/2023 18:45",6.19,"25/11/2023 19:13",5.89,"30/11/2023 18:51",8.77,"25/11/2023 19:13",9.36,"30/11/2023 18:51","https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ittihad-al-khaleej/fThVX3B9/")
)

$startRow = 128
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $targetRow = $startRow + $i
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value2 = $rowData[$c]
    }
}

Write-Host "edit complete"
